$d = $word.ActiveDocument

# The document ends with a list paragraph (ListParagraph style, numbered list,
# bottom border, justified) whose last run is a single "." character.
# We need to append a brand-new list paragraph with the same paragraph
# formatting after it, containing the new stipulation text.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range

# Isolate just the trailing "." character (a non-bold run) so that the
# paragraph-mark / run formatting picked up for the newly inserted paragraph
# is plain (not bold), matching the formatting of the rest of the document.
$dotStart = $lastRange.End - 2
$dotEnd = $lastRange.End - 1
$dotRange = $d.Range($dotStart, $dotEnd)

$find = $dotRange.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$newText = "Licensing or leasing or purchasing does not allow the payer who funded or conducted the illegal activities that are now blocked, to change the source code or systems, in any manner that would change them or make them not work."

$find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, ".^p" + $newText, 2) | Out-Null
